$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-cell value updates (D and/or E columns) ---
# Each entry: row, column letter, new value
$cellUpdates = @(
    @{ Row = 2;  Col = "D"; Value = "73.031.61" },
    @{ Row = 2;  Col = "E"; Value = "  +1.20%  " },

    @{ Row = 3;  Col = "D"; Value = "3.963.68" },
    @{ Row = 3;  Col = "E"; Value = "  -1.82%  " },

    @{ Row = 4;  Col = "E"; Value = "  -0.07%  " },

    @{ Row = 5;  Col = "D"; Value = "616.23" },
    @{ Row = 5;  Col = "E"; Value = "  +14.19%  " },

    @{ Row = 6;  Col = "D"; Value = "166.06" },
    @{ Row = 6;  Col = "E"; Value = "  +9.04%  " },

    @{ Row = 7;  Col = "D"; Value = "0.680" },
    @{ Row = 7;  Col = "E"; Value = "  -2.72%  " },

    @{ Row = 8;  Col = "D"; Value = "0.999" },
    @{ Row = 8;  Col = "E"; Value = "  +0.05%  " },

    @{ Row = 9;  Col = "D"; Value = "0.755" },
    @{ Row = 9;  Col = "E"; Value = "  +0.15%  " },

    @{ Row = 10; Col = "D"; Value = "0.183" },
    @{ Row = 10; Col = "E"; Value = "  +6.69%  " },

    @{ Row = 11; Col = "D"; Value = "56.05" },
    @{ Row = 11; Col = "E"; Value = "  +4.93%  " },

    @{ Row = 12; Col = "D"; Value = "0.0000334" },
    @{ Row = 12; Col = "E"; Value = "  +1.56%  " },

    @{ Row = 13; Col = "D"; Value = "11.13" },
    @{ Row = 13; Col = "E"; Value = "  +2.24%  " },

    @{ Row = 14; Col = "D"; Value = "4.594.76" },
    @{ Row = 14; Col = "E"; Value = "  -1.83%  " },

    @{ Row = 15; Col = "D"; Value = "3.973.56" },
    @{ Row = 15; Col = "E"; Value = "  -2.02%  " },

    @{ Row = 16; Col = "E"; Value = "  +3.71%  " },

    @{ Row = 17; Col = "D"; Value = "14.11" },
    @{ Row = 17; Col = "E"; Value = "  -1.87%  " },

    @{ Row = 18; Col = "D"; Value = "20.47" },
    @{ Row = 18; Col = "E"; Value = "  -0.63%  " },

    @{ Row = 21; Col = "D"; Value = "439.41" },
    @{ Row = 21; Col = "E"; Value = "  -1.92%  " },

    @{ Row = 22; Col = "D"; Value = "4.90" },
    @{ Row = 22; Col = "E"; Value = "  +15.33%  " },

    @{ Row = 23; Col = "D"; Value = "95.85" },
    @{ Row = 23; Col = "E"; Value = "  -1.83%  " },

    @{ Row = 24; Col = "D"; Value = "3.37" },
    @{ Row = 24; Col = "E"; Value = "  -3.85%  " },

    @{ Row = 25; Col = "D"; Value = "14.15" },
    @{ Row = 25; Col = "E"; Value = "  -3.09%  " },

    @{ Row = 26; Col = "D"; Value = "4.09" },
    @{ Row = 26; Col = "E"; Value = "  -6.14%  " },

    @{ Row = 27; Col = "D"; Value = "11.05" },
    @{ Row = 27; Col = "E"; Value = "  -1.97%  " },

    @{ Row = 28; Col = "D"; Value = "5.96" },
    @{ Row = 28; Col = "E"; Value = "  +0.08%  " },

    @{ Row = 29; Col = "D"; Value = "10.49" },
    @{ Row = 29; Col = "E"; Value = "  -2.74%  " },

    @{ Row = 30; Col = "D"; Value = "36.03" },
    @{ Row = 30; Col = "E"; Value = "  -2.97%  " },

    @{ Row = 31; Col = "E"; Value = "  -1.80%  " },

    @{ Row = 32; Col = "D"; Value = "13.62" },
    @{ Row = 32; Col = "E"; Value = "  +0.31%  " },

    @{ Row = 33; Col = "D"; Value = "0.0000105" },
    @{ Row = 33; Col = "E"; Value = "  +19.80%  " },

    @{ Row = 34; Col = "E"; Value = "  -3.84%  " },

    @{ Row = 35; Col = "D"; Value = "47.42" },
    @{ Row = 35; Col = "E"; Value = "  -3.57%  " },

    @{ Row = 36; Col = "D"; Value = "70.64" },
    @{ Row = 36; Col = "E"; Value = "  +5.64%  " },

    @{ Row = 37; Col = "D"; Value = "648.66" },
    @{ Row = 37; Col = "E"; Value = "  -4.51%  " },

    @{ Row = 38; Col = "E"; Value = "  -5.53%  " },

    @{ Row = 39; Col = "E"; Value = "  +2.40%  " },

    @{ Row = 40; Col = "E"; Value = "  -0.05%  " },

    @{ Row = 41; Col = "E"; Value = "  -1.85%  " },

    @{ Row = 42; Col = "E"; Value = "  +0.14%  " },

    @{ Row = 43; Col = "E"; Value = "  -2.26%  " },

    @{ Row = 44; Col = "D"; Value = "10.68" },
    @{ Row = 44; Col = "E"; Value = "  -4.85%  " },

    @{ Row = 45; Col = "E"; Value = "  -5.40%  " },

    @{ Row = 46; Col = "E"; Value = "  -1.47%  " },

    @{ Row = 47; Col = "E"; Value = "  +3.30%  " },

    @{ Row = 48; Col = "D"; Value = "2.90" },
    @{ Row = 48; Col = "E"; Value = "  +26.82%  " },

    @{ Row = 49; Col = "D"; Value = "2.59" },
    @{ Row = 49; Col = "E"; Value = "  -1.54%  " }
)

# Price figures in column D look numeric (e.g. "616.23") and Excel would
# silently coerce them to real numbers (losing trailing zeros, switching
# to scientific notation for tiny values, etc.) unless we tell it to treat
# the entry as text. Prefixing with a leading apostrophe mirrors typing the
# value into the Excel UI as text and keeps the original "General" cell
# formatting (no NumberFormat/style changes), only marking quotePrefix.
foreach ($u in $cellUpdates) {
    $value = $u.Value
    if ($u.Col -eq "D") {
        $value = "'" + $value
    }
    $ws.Range("$($u.Col)$($u.Row)").Value = $value
}

# --- Full row swaps (Coin/Link/Price/Volume) ---
# Row 19 used to hold WrappedBTC, row 20 held TRON; they swap places
# and each gets refreshed Price/Volume figures.
$rowUpdates = @(
    @{ Row = 19; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.131"; E = "  -0.41%  " },
    @{ Row = 20; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "72.832.06"; E = "  +0.98%  " },
    # Row 50 used to hold Maker, row 51 held Monero; Maker moves down to
    # row 51 and row 50 now lists FLOKI.
    @{ Row = 50; B = "FLOKI"; C = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"; D = "0.000280"; E = "  +1.28%  " },
    @{ Row = 51; B = "Maker"; C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D = "2.832.38"; E = "  +3.35%  " }
)

foreach ($r in $rowUpdates) {
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = "'" + $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
}
